$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts old D:K data to F:M)
$ws.Range("D:E").Insert()

# Step 2: Copy cell formatting (number formats/styles) from the shifted data (now in F:G)
# into the newly inserted D:E columns, for the full data block (rows 5-102).
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the column width of the new columns to the rest of the quarterly data columns.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# Step 3: Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with values.

# 3a. Numeric values (including dates stored as serials, and zeros)
$numericRows = @{
    7 = @(43465,43373)
    8 = @(1864000,1634000)
    13 = @(0,0)
    14 = @(0,0)
    15 = @(-35000,-48000)
    17 = @(979000,740000)
    18 = @(885000,894000)
    20 = @(135000,440000)
    21 = @(1348000,1668000)
    22 = @(0,0)
    23 = @(1020000,1334000)
    24 = @(163000,313000)
    25 = @(0,0)
    26 = @(857000,1021000)
    27 = @(814000,975000)
    28 = @(0,0)
    29 = @(13000,93000)
    30 = @(0,0)
    31 = @(0,0)
    32 = @(-135000,-440000)
    33 = @(827000,1068000)
    34 = @(0,0)
    35 = @(827000,1068000)
    38 = @(43465,43373)
    41 = @(73852000,79772000)
    42 = @(68730000,51891000)
    43 = @(0,0)
    44 = @(0,0)
    45 = @(0,0)
    46 = @(0,0)
    47 = @(1064000,1143000)
    48 = @(1832000,1832000)
    49 = @(22222000,22227000)
    50 = @(0,0)
    51 = @(0,0)
    52 = @(1357000,1558000)
    53 = @(0,0)
    54 = @(362873000,349770000)
    57 = @(19731000,18683000)
    58 = @(0,0)
    59 = @(5669000,5601000)
    60 = @(0,0)
    61 = @(29163000,28113000)
    62 = @(0,0)
    63 = @(0,0)
    64 = @(0,0)
    65 = @(0,0)
    66 = @(322235000,308210000)
    68 = @(0,0)
    69 = @(0,0)
    70 = @(3542000,3542000)
    71 = @(0,0)
    72 = @(28652000,28098000)
    73 = @(0,0)
    74 = @(0,0)
    75 = @(0,0)
    76 = @(37096000,38018000)
    77 = @(0,0)
    80 = @(43465,43373)
    81 = @(827000,1068000)
    83 = @(328000,334000)
    84 = @(0,0)
    85 = @(0,0)
    86 = @(0,0)
    87 = @(0,0)
    88 = @(0,0)
    89 = @(3184000,1244000)
    91 = @(-289000,-314000)
    92 = @(0,0)
    93 = @(0,0)
    94 = @(-14718000,2803000)
    96 = @(-278000,-283000)
    97 = @(0,0)
    98 = @(0,0)
    99 = @(0,0)
    100 = @(13458000,-4721000)
    101 = @(-15000,3000)
    102 = @(1909000,-671000)
}
foreach ($r in $numericRows.Keys) {
    $vals = $numericRows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("E$r").Value = $vals[1]
}

# 3b. Rows that display the "NA" placeholder text in the new columns
$naRows = @(9,10,12)
foreach ($r in $naRows) {
    $ws.Range("D$r").Value = "NA"
    $ws.Range("E$r").Value = "NA"
}

# 3c. Rows that stay blank in the new columns (spacer rows) -- nothing to do,
#     they were left blank by the format-only paste above.

# Step 4: A handful of historical quarterly figures were also corrected as part of
#         this data refresh (same cells identified in the published diff).
$ws.Range("I89").Value = 5212000
$ws.Range("I91").Value = -211000
$ws.Range("J91").Value = -436000
$ws.Range("H94").Value = -17549000
$ws.Range("I94").Value = -73000
$ws.Range("H102").Value = -701000
$ws.Range("I102").Value = 1358000

